# Update the parameter-estimation config start_values (column B) so the
# workbook works with COPASI version 19 and the PE report gives proper
# titles (per commit message). Only the "start_value" column changes for
# four parameters:
#   (reaction_3).k1        : 0.1                  -> 10
#   assignment_global_var  : 1.999999999999998     -> 1.999999999999999
#   two                    : 50                    -> 15
#   B (species)            : 9.99999845072         -> 5.0
#
# These replacement values look numeric, but the column stores them as
# text (shared strings) in the workbook, matching every other value in
# that column (e.g. "0.1", "1e-06", ...). Assigning a numeric-looking
# string straight to Range.Value would make Excel coerce it to a real
# number, so each cell is briefly forced to Text format, written, then
# restored to the default "Normal" style (clearing the explicit
# NumberFormat again) so the cell keeps its original, unstyled look.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param(
        [string]$CellRef,
        [string]$Text
    )
    $rng = $ws.Range($CellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $Text
    $rng.Style = "Normal"
}

Set-TextValue "B2"  "10"
Set-TextValue "B7"  "1.999999999999999"
Set-TextValue "B8"  "15"
Set-TextValue "B12" "5.0"
